# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Map of row -> new value for column G (header "K")
$newK = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 2
    23 = 2
    24 = 2
    25 = 0
    26 = 1
    27 = 3
    28 = 2
    29 = 1
    30 = 2
    31 = 1
    33 = 2
    34 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
